$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 previously held the text "R40"; update it to hold the text "1".
# A leading apostrophe forces Excel to store a numeric-looking entry as text
# (shared string) instead of auto-converting it to a number.
$ws.Range("B11").Value = "'1"
